$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 currently holds the placeholder "xxx" -> rename to the real hotel name
$ws.Range("B2").Value = "Jumeirah Beach Hotel"

# B3 gets a new, distinct hotel name (currently duplicates B2's old value)
$ws.Range("B3").Value = "Grand Plaza Apartments"
